# Final changes 17th March 2022
# Update ShipmentTracking (col P) numbers for rows 2-25 and the
# ActualRate (col Q) value for row 24. Values are entered with a
# leading apostrophe so Excel stores them as text (matching the
# existing t="s" shared-string cell type) instead of auto-converting
# the digit-only / currency-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Value  = "'320018063820"
$ws.Range("P3").Value  = "'320018063831"
$ws.Range("P4").Value  = "'320018063864"
$ws.Range("P5").Value  = "'320018063886"
$ws.Range("P6").Value  = "'320018063923"
$ws.Range("P7").Value  = "'320018063945"
$ws.Range("P8").Value  = "'320018063978"
$ws.Range("P9").Value  = "'320018063990"
$ws.Range("P10").Value = "'320018064025"
$ws.Range("P11").Value = "'320018064047"
$ws.Range("P12").Value = "'320018064080"
$ws.Range("P13").Value = "'320018064106"
$ws.Range("P14").Value = "'320018064139"
$ws.Range("P15").Value = "'320018064150"
$ws.Range("P16").Value = "'320018064183"
$ws.Range("P17").Value = "'320018064209"
$ws.Range("P18").Value = "'320018064242"
$ws.Range("P19").Value = "'320018064264"
$ws.Range("P20").Value = "'320018064297"
$ws.Range("P21").Value = "'320018064312"
$ws.Range("P22").Value = "'320018064345"
$ws.Range("P23").Value = "'320018064356"
$ws.Range("P24").Value = "'320018064367"
$ws.Range("Q24").Value = "'$248.51"
$ws.Range("P25").Value = "'320018064378"
